$d = $word.ActiveDocument

# Locate the "Research & Data Analytics Leadership" paragraph under the Siege Analytics
# (PARTNER) role, immediately before the existing bullet list.
$rng = $d.Content
$rng.Find.Execute("Research & Data Analytics Leadership", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)

$newBullets = @(
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

foreach ($bulletText in $newBullets) {
    $rng.InsertParagraphAfter() | Out-Null
    $rng.Move(4, 1) | Out-Null
    $rng.InsertBefore($bulletText)
    $rng.Collapse(0)
}

$d.Save()
